$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "57.421.82"
$ws.Range("E2").Value = "  +0.28%  "
$ws.Range("D3").Value = "2.414.24"
$ws.Range("E3").Value = "  -4.18%  "
$ws.Range("E4").Value = "  -0.15%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "488.32"
$ws.Range("D5").ClearFormats()
$ws.Range("E5").Value = "  -1.65%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "153.35"
$ws.Range("D6").ClearFormats()
$ws.Range("E6").Value = "  -0.33%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.997"
$ws.Range("D7").ClearFormats()
$ws.Range("E7").Value = "  +0.14%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.612"
$ws.Range("D8").ClearFormats()
$ws.Range("E8").Value = "  +18.02%  "
$ws.Range("D9").Value = "2.434.28"
$ws.Range("E9").Value = "  -4.07%  "
$ws.Range("E10").Value = "  +5.61%  "
$ws.Range("E11").Value = "  -1.65%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.332"
$ws.Range("D12").ClearFormats()
$ws.Range("E12").Value = "  -2.27%  "
$ws.Range("E13").Value = "  +1.18%  "
$ws.Range("D14").Value = "2.839.76"
$ws.Range("E14").Value = "  -4.04%  "
$ws.Range("D15").Value = "57.354.20"
$ws.Range("E15").Value = "  -0.04%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "20.55"
$ws.Range("D16").ClearFormats()
$ws.Range("E16").Value = "  -4.27%  "
$ws.Range("E17").Value = "  -4.33%  "
$ws.Range("D18").Value = "2.435.50"
$ws.Range("E18").Value = "  -3.49%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "4.65"
$ws.Range("D19").ClearFormats()
$ws.Range("E19").Value = "  +1.05%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "323.74"
$ws.Range("D20").ClearFormats()
$ws.Range("E20").Value = "  -0.39%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "9.99"
$ws.Range("D21").ClearFormats()
$ws.Range("E21").Value = "  -3.79%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "1.00"
$ws.Range("D22").ClearFormats()
$ws.Range("E22").Value = "  +0.14%  "
$ws.Range("E23").Value = "  -0.90%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "57.93"
$ws.Range("D24").ClearFormats()
$ws.Range("E24").Value = "  -1.26%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "0.403"
$ws.Range("D25").ClearFormats()
$ws.Range("E25").Value = "  -2.31%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "0.996"
$ws.Range("D26").ClearFormats()
$ws.Range("E26").Value = "  -0.07%  "
$ws.Range("E27").Value = "  -3.52%  "
$ws.Range("D28").Value = "2.522.83"
$ws.Range("E28").Value = "  -3.75%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "7.27"
$ws.Range("D29").ClearFormats()
$ws.Range("E29").Value = "  -5.14%  "
$ws.Range("E30").Value = "  -6.00%  "
$ws.Range("E31").Value = "  +0.16%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "150.92"
$ws.Range("D32").ClearFormats()
$ws.Range("E32").Value = "  -0.66%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "18.62"
$ws.Range("D33").ClearFormats()
$ws.Range("E33").Value = "  +0.88%  "
$ws.Range("E34").Value = "  -1.52%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "5.29"
$ws.Range("D35").ClearFormats()
$ws.Range("E35").Value = "  -0.22%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "1.15"
$ws.Range("D36").ClearFormats()
$ws.Range("E36").Value = "  -1.55%  "
$ws.Range("E37").Value = "  -2.71%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.816"
$ws.Range("D38").ClearFormats()
$ws.Range("E38").Value = "  -11.30%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.101"
$ws.Range("D39").ClearFormats()
$ws.Range("E39").Value = "  +6.27%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "33.97"
$ws.Range("D40").ClearFormats()
$ws.Range("E40").Value = "  -1.39%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "283.56"
$ws.Range("D41").ClearFormats()
$ws.Range("E41").Value = "  +4.96%  "
$ws.Range("E42").Value = "  -1.16%  "
$ws.Range("E43").Value = "  -3.46%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.997"
$ws.Range("D44").ClearFormats()
$ws.Range("E44").Value = "  +0.21%  "
$ws.Range("E45").Value = "  -3.81%  "
$ws.Range("E46").Value = "  -6.26%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "10.24"
$ws.Range("D47").ClearFormats()
$ws.Range("E47").Value = "  +0.17%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "0.0227"
$ws.Range("D48").ClearFormats()
$ws.Range("E48").Value = "  -1.64%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "4.53"
$ws.Range("D49").ClearFormats()
$ws.Range("E49").Value = "  -7.50%  "
$ws.Range("D50").Value = "1.899.94"
$ws.Range("E50").Value = "  -0.35%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "17.60"
$ws.Range("D51").ClearFormats()
$ws.Range("E51").Value = "  -3.53%  "
